$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 headers only: insert "Bucket Location" before "File Tag 1",
# shift the existing "File Tag 1"/"File Tag 2" headers right, and add
# "File Tag 3"/"File Tag 4"/"File Tag 5" — without touching the data
# rows below (they already have values in columns C:H).
$ws.Range("C1").Value = "Bucket Location"
$ws.Range("D1").Value = "File Tag 1"
$ws.Range("E1").Value = "File Tag 2"
$ws.Range("F1").Value = "File Tag 3"
$ws.Range("G1").Value = "File Tag 4"
$ws.Range("H1").Value = "File Tag 5"
